$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 71; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value()
    if ($v -eq "Masculino") {
        $cell.Value = "Hombre"
    } elseif ($v -eq "Femenino") {
        $cell.Value = "Mujer"
    }
}

$ws.Range("R60").Select()
